$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2024291497975708
$ws.Range("C2").Value = 0.5668016194331984
$ws.Range("J2").Value = 0.008097165991902834
$ws.Range("P2").Value = 0.1295546558704453
$ws.Range("S2").Value = 0.0931174089068826
$ws.Range("B3").Value = 0.01398601398601399
$ws.Range("C3").Value = 0.01398601398601399
$ws.Range("J3").Value = 0.04195804195804196
$ws.Range("P3").Value = 0.7832167832167832
$ws.Range("S3").Value = 0.1468531468531468
$ws.Range("P4").Value = 0.7857142857142857
$ws.Range("S4").Value = 0.2142857142857143
$ws.Range("B6").Value = 0.04145077720207254
$ws.Range("D6").Value = 0.01036269430051814
$ws.Range("F6").Value = 0.06217616580310881
$ws.Range("J6").Value = 0.3056994818652849
$ws.Range("O6").Value = 0.02072538860103627
$ws.Range("Q6").Value = 0.150259067357513
$ws.Range("R6").Value = 0.09844559585492228
$ws.Range("S6").Value = 0.310880829015544
$ws.Range("B7").Value = 0.1125
$ws.Range("D7").Value = 0.00625
$ws.Range("E7").Value = 0.00625
$ws.Range("F7").Value = 0.0375
$ws.Range("J7").Value = 0.14375
$ws.Range("O7").Value = 0.00625
$ws.Range("Q7").Value = 0.19375
$ws.Range("R7").Value = 0.06875000000000001
$ws.Range("S7").Value = 0.425
$ws.Range("B8").Value = 0.09424083769633508
$ws.Range("D8").Value = 0.01047120418848168
$ws.Range("F8").Value = 0.03926701570680628
$ws.Range("J8").Value = 0.1047120418848168
$ws.Range("O8").Value = 0.02094240837696335
$ws.Range("Q8").Value = 0.2356020942408377
$ws.Range("R8").Value = 0.09424083769633508
$ws.Range("S8").Value = 0.4005235602094241
$ws.Range("D9").Value = 0.01515151515151515
$ws.Range("F9").Value = 0.04545454545454546
$ws.Range("J9").Value = 0.1363636363636364
$ws.Range("O9").Value = 0.007575757575757576
$ws.Range("Q9").Value = 0.2121212121212121
$ws.Range("R9").Value = 0.1136363636363636
$ws.Range("S9").Value = 0.3863636363636364
$ws.Range("B10").Value = 0.09256449165402124
$ws.Range("D10").Value = 0.02655538694992413
$ws.Range("E10").Value = 0.001517450682852807
$ws.Range("F10").Value = 0.0629742033383915
$ws.Range("J10").Value = 0.1358118361153262
$ws.Range("O10").Value = 0.009863429438543247
$ws.Range("Q10").Value = 0.2397572078907436
$ws.Range("R10").Value = 0.07587253414264036
$ws.Range("S10").Value = 0.3550834597875569
$ws.Range("G11").Value = 0.1830985915492958
$ws.Range("J11").Value = 0.102112676056338
$ws.Range("K11").Value = 0.2323943661971831
$ws.Range("L11").Value = 0.4436619718309859
$ws.Range("S11").Value = 0.03873239436619718
$ws.Range("G12").Value = 0.7022900763358778
$ws.Range("J12").Value = 0.2290076335877863
$ws.Range("K12").Value = 0.007633587786259542
$ws.Range("L12").Value = 0.01526717557251908
$ws.Range("S12").Value = 0.04580152671755725
$ws.Range("G13").Value = 0.6129032258064516
$ws.Range("J13").Value = 0.3548387096774194
$ws.Range("S13").Value = 0.03225806451612903
$ws.Range("F15").Value = 0.009049773755656109
$ws.Range("H15").Value = 0.1040723981900453
$ws.Range("I15").Value = 0.06334841628959276
$ws.Range("J15").Value = 0.4570135746606335
$ws.Range("K15").Value = 0.06334841628959276
$ws.Range("M15").Value = 0.004524886877828055
$ws.Range("O15").Value = 0.05429864253393665
$ws.Range("S15").Value = 0.2443438914027149
$ws.Range("F16").Value = 0.01169590643274854
$ws.Range("H16").Value = 0.1345029239766082
$ws.Range("I16").Value = 0.07017543859649122
$ws.Range("J16").Value = 0.5029239766081871
$ws.Range("K16").Value = 0.08187134502923976
$ws.Range("M16").Value = 0.01169590643274854
$ws.Range("O16").Value = 0.08187134502923976
$ws.Range("S16").Value = 0.1052631578947368
$ws.Range("F17").Value = 0.03258655804480651
$ws.Range("H17").Value = 0.1812627291242362
$ws.Range("I17").Value = 0.05906313645621181
$ws.Range("J17").Value = 0.4215885947046843
$ws.Range("K17").Value = 0.1038696537678208
$ws.Range("M17").Value = 0.01629327902240326
$ws.Range("O17").Value = 0.06313645621181263
$ws.Range("S17").Value = 0.1221995926680244
$ws.Range("F18").Value = 0.03333333333333333
$ws.Range("H18").Value = 0.1888888888888889
$ws.Range("I18").Value = 0.06111111111111111
$ws.Range("J18").Value = 0.4444444444444444
$ws.Range("K18").Value = 0.1
$ws.Range("M18").Value = 0.01666666666666667
$ws.Range("O18").Value = 0.07222222222222222
$ws.Range("S18").Value = 0.08333333333333333
$ws.Range("F19").Value = 0.01898188093183779
$ws.Range("H19").Value = 0.1872303710094909
$ws.Range("I19").Value = 0.05694564279551338
$ws.Range("J19").Value = 0.4012079378774806
$ws.Range("K19").Value = 0.09836065573770492
$ws.Range("M19").Value = 0.01639344262295082
$ws.Range("O19").Value = 0.08369283865401209
$ws.Range("S19").Value = 0.1371872303710095
